# Insert a new data row at row 27 (pushes existing rows 27-52 down to 28-53)
# and populate it with the new "Ají" record described in the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(27).Insert()

$ws.Cells.Item(27, 1).Value = 7
$ws.Cells.Item(27, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(27, 3).Value = "Ñuble"
$ws.Cells.Item(27, 4).Value = 44566
$ws.Cells.Item(27, 4).NumberFormat = $ws.Cells.Item(28, 4).NumberFormat
$ws.Cells.Item(27, 5).Value = 16
$ws.Cells.Item(27, 6).Value = 100112021
$ws.Cells.Item(27, 7).Value = "Ají"
$ws.Cells.Item(27, 8).Value = "Americana (o)"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 60
$ws.Cells.Item(27, 11).Value = 20000
$ws.Cells.Item(27, 12).Value = 21000
$ws.Cells.Item(27, 13).Value = 20500
$ws.Cells.Item(27, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(27, 15).Value = "Región del Maule"
$ws.Cells.Item(27, 16).Value = 1367
$ws.Cells.Item(27, 17).Value = 15
$ws.Cells.Item(27, 18).Value = "Hortaliza"
